# Week 7_13.React JS-HOL_hands on.docx — title line correction:
#   "Week 6-13" + ".React JS-HOL"  ->  "Week 7" + "-13.React JS-HOL"
# with the stray "_GoBack" bookmark relocated from the end of the
# document to the split point between the two title runs (this mirrors
# Word's own behaviour of keeping only one "_GoBack" bookmark and moving
# it to the most-recent edit position).

$d = $word.ActiveDocument

# Locate the two runs that make up the title line.
$weekRange = $d.Content.Duplicate
$weekRange.Find.Execute("Week 6-13")

$reactRange = $d.Content.Duplicate
$reactRange.Find.Execute(".React JS-HOL")

# Drop a bookmark at the boundary between the two runs *before* editing
# any text - Word only ever keeps a single "_GoBack" bookmark, so adding
# it here automatically removes the old one near the end of the body.
$splitPoint = $d.Range($weekRange.End, $weekRange.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# Edit the second (right-hand) run first so the still-unmodified first
# run's cached Start/End positions stay valid.
$reactRange.Text = "-13.React JS-HOL"
$weekRange.Text = "Week 7"

Write-Output "done"
